$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

$wsZh.Range("D5").Value = "2016-01-18 03:35:32"
$wsDe.Range("D5").Value = "2016-01-18 03:35:47"
